$wb = $excel.ActiveWorkbook

# --- Sheet2: move the view's selection to a full-column range (A1:A16) ---
# Select the range on Sheet2 first (while Sheet1 is still the active/tabbed
# sheet) so Sheet2 keeps a stored selection of A1:A16 without becoming the
# active tab.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A1:A16").Select()

# --- Sheet1: becomes the active/selected tab ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

# Add the new value in A6
$ws1.Range("A6").Value = 12345

# Move Sheet1's selection to B6 (was Q15)
$ws1.Range("B6").Select()
